# Append two new log rows (153 and 154, 1-indexed data rows 154/155 in the
# sheet) to the feed_logs worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A154").Value = 153
$ws.Range("B154").Value = 1
$ws.Range("C154").Value = "2024-06-18 05:13:17"
$ws.Range("D154").Value = 200
$ws.Range("E154").Value = 10

$ws.Range("A155").Value = 154
$ws.Range("B155").Value = 2
$ws.Range("C155").Value = "2024-06-18 05:13:17"
$ws.Range("D155").Value = 200
$ws.Range("E155").Value = 1
